# Append six new list items to the end of the document, matching the
# existing "ListParagraph" bullet style used throughout the list.

$d = $word.ActiveDocument

$newLines = @(
    "Git checkout -b creates a new branch ",
    "Seems like you need to commit something for the git branch command to show the star",
    "Git diff nameofbranch will show the differences ",
    "Can push changes from a feature branch up to github. Git asks for an upstream branch and suggests something. I can than push with that nusach which they give me. ",
    "So, you can push from the feature branch up to remote branch. Merge there and than pull down to the local main branch. ",
    "Git branch -d nameofbranch deletes the branch which you will do after merging it in cuz why do you need it "
)

foreach ($line in $newLines) {
    $para = $d.Paragraphs.Last
    $rng = $para.Range
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $rng.InsertAfter($line)
}
